$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.317.42"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.867.77"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2802"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").Value = "1.862.41"
$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07457"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.070"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6535"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").Value = "30.278.29"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007582"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("D20").Value = "2.105.80"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.270"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.289"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.968"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09328"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.295"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.004"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05038"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.205"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7438"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.719"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01827"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9117"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.075"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.910"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4259"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.383"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1289"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.472"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05635"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.97%  "

